$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025-2")

# Row 13 - CHATA EXALMAR
$ws.Range("E13").Value = "A.S/0047-225"
$ws.Range("D13").Value = "A.S/0047"
$ws.Range("B13").Value = "CHATA EXALMAR"
$ws.Range("C13").Value = "Chata"
$ws.Range("A13").Value = "2025-2"

# Row 14 - EP MODESTO 7
$ws.Range("D14").Value = "A.S/0049"
$ws.Range("E14").Value = "A.S/0049-225"
$ws.Range("B14").Value = "EP MODESTO 7"
$ws.Range("C14").Value = "Embarcación Pesquera"
$ws.Range("A14").Value = "2025-2"

# Row 15 - EP MODESTO 9
$ws.Range("B15").Value = "EP MODESTO 9"
$ws.Range("D15").Value = "A.S/0050"
$ws.Range("E15").Value = "A.S/0050-225"
$ws.Range("C15").Value = "Embarcación Pesquera"
$ws.Range("A15").Value = "2025-2"

# Row 16 - EP DON MILTON
$ws.Range("B16").Value = "EP DON MILTON"
$ws.Range("D16").Value = "A.S/0048"
$ws.Range("E16").Value = "A.S/0048-225"
$ws.Range("C16").Value = "Embarcación Pesquera"
$ws.Range("A16").Value = "2025-2"

# Row 17 - EP DANIA
$ws.Range("B17").Value = "EP DANIA"
$ws.Range("D17").Value = "A.S/0051"
$ws.Range("E17").Value = "A.S/0051-225"
$ws.Range("C17").Value = "Embarcación Pesquera"
$ws.Range("A17").Value = "2025-2"

# Apply the same style as other cells in column C (centered alignment, matching font/border)
$ws.Range("C12").Copy()
$ws.Range("C13:C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column width adjustments to match auto-fit after the new, wider entries
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 21
$ws.Columns.Item(5).ColumnWidth = 13.33

# Update selection to match final state left by the author
$ws.Range("E20").Select()
